# "change config directory back to red cross demo, also edit forms to be
#  compatable with rev 210 translations changes"
#
# - Makes the "survey" sheet the active tab again (was "properties").
# - Updates the "survey" sheet's header cell C1 from "display.text" to the
#   new "display.prompt.text" label (rev 210 translations change), which
#   also selects C1 (instead of C2) as the active cell on that sheet.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$survey.Range("C1").Value = "display.prompt.text"
$survey.Range("C1").Select()
$survey.Activate()
